$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Sarfaraz Khan"

# Numeric-looking text columns (runs, balls, fours, sixes, sr) must stay text,
# like the rest of the sheet, so force text format before writing values.
$ws.Range("E2:I3").NumberFormat = "@"

# Header row (row 1) - insert matchNo as new column A, shift rest right
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# Row 2 - new match record (53rd match vs CSK)
$ws.Range("A2").Value = "53rd"
$ws.Range("B2").Value = "Punjab Kings"
$ws.Range("C2").Value = "Sarfaraz Khan"
$ws.Range("D2").Value = "c du Plessis b Thakur"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "3"
$ws.Range("G2").Value = "0"
$ws.Range("H2").Value = "0"
$ws.Range("I2").Value = "0.00"
$ws.Range("J2").Value = "Chennai Super Kings"
$ws.Range("K2").Value = "Dubai (DSC)"
$ws.Range("L2").Value = "October 07"
$ws.Range("M2").Value = "Punjab Kings won by 6 wickets (with 42 balls remaining)"

# Row 3 - original match record (48th match vs RCB), shifted right by one column
$ws.Range("A3").Value = "48th"
$ws.Range("B3").Value = "Punjab Kings"
$ws.Range("C3").Value = "Sarfaraz Khan"
$ws.Range("D3").Value = "b Chahal"
$ws.Range("E3").Value = "0"
$ws.Range("F3").Value = "1"
$ws.Range("G3").Value = "0"
$ws.Range("H3").Value = "0"
$ws.Range("I3").Value = "0.00"
$ws.Range("J3").Value = "Royal Challengers Bangalore"
$ws.Range("K3").Value = "Sharjah"
$ws.Range("L3").Value = "October 03"
$ws.Range("M3").Value = "RCB won by 6 runs"
